# Add a "mean +/- SD" summary table on a new worksheet.
#
# 1. Rename the original sheet to "Overall LIC".
# 2. Duplicate it to create "Spatially Explicit LIC" right after it (this
#    keeps the existing fonts/styles/row heights instead of minting new
#    style entries) and make it the active sheet/tab.
# 3. Drop the duplicated H:L "spatially explicit" block (columns H:L are
#    no longer needed once this becomes its own sheet) and replace the
#    A:E values with the Species / decade header plus the mean +/- SD
#    values.
# 4. Leave the "Overall LIC" sheet's selection on A1:E7 (its own data
#    block) since it is no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- sheet1: rename -------------------------------------------------
$overall = $wb.Worksheets.Item(1)
$overall.Name = "Overall LIC"

# --- new sheet2: duplicate "Overall LIC" right after itself ---------
$overall.Copy($null, $overall)
$spatial = $wb.Worksheets.Item(2)
$spatial.Name = "Spatially Explicit LIC"

# Remove the duplicated H:L block - the new sheet only needs A:E.
$spatial.Range("H1:L7").Delete()

# --- header row -------------------------------------------------------
$spatial.Range("A1").Value = "Species"
$spatial.Range("B1").Value = "1980s"
$spatial.Range("C1").Value = "1990s"
$spatial.Range("D1").Value = "2000s"
$spatial.Range("E1").Value = "2010s"

# --- data rows (mean +/- SD), replacing the copied numeric values -----
$data = @(
    @("starry flounder",   "0.0006 ± 0.004",  "0.00005 ± 0.001", "0.0009 ± 0.006", "0.0002 ± 0.001"),
    @("petrale sole",      "0.0009 ± 0.003",  "0.0033 ± 0.006",  "0.0057 ± 0.008", "0.0049 ± 0.006"),
    @("sand sole",         "0.0012 ± 0.007",  "0.0013 ± 0.011",  "0.0022 ± 0.013", "0.0001 ± 0.001"),
    @("English sole",      "0.0025 ± 0.006",  "0.0020 ± 0.004",  "0.0038 ± 0.007", "0.0032 ± 0.013"),
    @("Pacific sanddab",   "0.0036 ± 0.024",  "0.0042 ± 0.017",  "0.0038 ± 0.013", "0.0058 ± 0.020"),
    @("Dover sole",        "0.0046 ± 0.010",  "0.0054 ± 0.012",  "0.0059 ± 0.011", "0.0046 ± 0.008")
)

$row = 2
foreach ($r in $data) {
    $spatial.Cells.Item($row, 1).Value = $r[0]
    $spatial.Cells.Item($row, 2).Value = $r[1]
    $spatial.Cells.Item($row, 3).Value = $r[2]
    $spatial.Cells.Item($row, 4).Value = $r[3]
    $spatial.Cells.Item($row, 5).Value = $r[4]
    # Single-line mean +/- SD text no longer needs the 2-line wrapped
    # row height inherited from the copied sheet.
    $spatial.Rows.Item($row).RowHeight = 15.5
    $row = $row + 1
}

# --- selections / active tab ------------------------------------------
$overall.Range("A1:E7").Select()
$spatial.Activate()
$spatial.Range("D13").Select()
